# Generate Report for Handback
#
# The localization-status workbook tracks handoff/handback state for each
# language. This run represents a handback event: the zh-cn and de-de
# files just came back from translation, so:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - The "Latest Target File" / "Latest Handback File" columns (E/F) get
#     populated with links to the source markdown and the returned xlf
#   - The "Latest Handback DateTime" column (G) is stamped with the time
#     the handback was recorded.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# --- Status column: flip every "Ready for handoff" cell to the handed-back
#     message. Doing all of them together lets the shared-string table
#     collapse back down instead of leaving an orphaned old string around.
$ws1.Range("B2").Value = $newStatus
$ws1.Range("C2").Value = $newStatus
$ws1.Range("B3").Value = $newStatus
$ws1.Range("C3").Value = $newStatus

$ws2.Range("B2").Value = $newStatus
$ws2.Range("B3").Value = $newStatus

$ws3.Range("B2").Value = $newStatus
$ws3.Range("B3").Value = $newStatus

# --- zh-cn sheet: populate Latest Target File (E) / Latest Handback File (F)
#     with hyperlinks mirroring the existing Source File Name (A) / Latest
#     Handoff File (C) links, and stamp Latest Handback DateTime (G).

$ws2.Range("E2").Value = "009578f8-0ead-4dfa-bb7f-6502baaccc9c.md"
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/e3ae61c96bf6b7d8d01ebc65d072ef585ab5a460/e2e/009578f8-0ead-4dfa-bb7f-6502baaccc9c.md", [Type]::Missing, [Type]::Missing, "009578f8-0ead-4dfa-bb7f-6502baaccc9c.md")

$ws2.Range("F2").Value = "009578f8-0ead-4dfa-bb7f-6502baaccc9c.8a33fa58a8effc5606d71c6558226c6904bfb48b.zh-cn.xlf"
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3622baa375f126b1102418b125be3cffd3c235d5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/009578f8-0ead-4dfa-bb7f-6502baaccc9c.8a33fa58a8effc5606d71c6558226c6904bfb48b.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "009578f8-0ead-4dfa-bb7f-6502baaccc9c.8a33fa58a8effc5606d71c6558226c6904bfb48b.zh-cn.xlf")

$ws2.Range("G2").Value = "2016-01-28 05:27:58"

$ws2.Range("E3").Value = "fa496268-6053-49df-8127-4f8d4c5800aa.md"
$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/e3ae61c96bf6b7d8d01ebc65d072ef585ab5a460/e2e/fa496268-6053-49df-8127-4f8d4c5800aa.md", [Type]::Missing, [Type]::Missing, "fa496268-6053-49df-8127-4f8d4c5800aa.md")

$ws2.Range("F3").Value = "fa496268-6053-49df-8127-4f8d4c5800aa.69f8b9e763b65ba9f489ceba2566a1f0bfdbabcb.zh-cn.xlf"
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3622baa375f126b1102418b125be3cffd3c235d5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/fa496268-6053-49df-8127-4f8d4c5800aa.69f8b9e763b65ba9f489ceba2566a1f0bfdbabcb.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "fa496268-6053-49df-8127-4f8d4c5800aa.69f8b9e763b65ba9f489ceba2566a1f0bfdbabcb.zh-cn.xlf")

$ws2.Range("G3").Value = "2016-01-28 05:27:58"

# --- de-de sheet: same pattern, its own handback timestamp.

$ws3.Range("E2").Value = "009578f8-0ead-4dfa-bb7f-6502baaccc9c.md"
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/e3ae61c96bf6b7d8d01ebc65d072ef585ab5a460/e2e/009578f8-0ead-4dfa-bb7f-6502baaccc9c.md", [Type]::Missing, [Type]::Missing, "009578f8-0ead-4dfa-bb7f-6502baaccc9c.md")

$ws3.Range("F2").Value = "009578f8-0ead-4dfa-bb7f-6502baaccc9c.8a33fa58a8effc5606d71c6558226c6904bfb48b.de-de.xlf"
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2b3a6d401e13231a16991b306922330fc247f4fa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/009578f8-0ead-4dfa-bb7f-6502baaccc9c.8a33fa58a8effc5606d71c6558226c6904bfb48b.de-de.xlf", [Type]::Missing, [Type]::Missing, "009578f8-0ead-4dfa-bb7f-6502baaccc9c.8a33fa58a8effc5606d71c6558226c6904bfb48b.de-de.xlf")

$ws3.Range("G2").Value = "2016-01-28 05:28:16"

$ws3.Range("E3").Value = "fa496268-6053-49df-8127-4f8d4c5800aa.md"
$ws3.Hyperlinks.Add($ws3.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/e3ae61c96bf6b7d8d01ebc65d072ef585ab5a460/e2e/fa496268-6053-49df-8127-4f8d4c5800aa.md", [Type]::Missing, [Type]::Missing, "fa496268-6053-49df-8127-4f8d4c5800aa.md")

$ws3.Range("F3").Value = "fa496268-6053-49df-8127-4f8d4c5800aa.69f8b9e763b65ba9f489ceba2566a1f0bfdbabcb.de-de.xlf"
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2b3a6d401e13231a16991b306922330fc247f4fa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/fa496268-6053-49df-8127-4f8d4c5800aa.69f8b9e763b65ba9f489ceba2566a1f0bfdbabcb.de-de.xlf", [Type]::Missing, [Type]::Missing, "fa496268-6053-49df-8127-4f8d4c5800aa.69f8b9e763b65ba9f489ceba2566a1f0bfdbabcb.de-de.xlf")

$ws3.Range("G3").Value = "2016-01-28 05:28:16"
